# Update "Estado de Cuenta" worker/period table on Hoja1.
# The original table listed each worker's arrears periods from the most
# recent (2410) down to the oldest (2404). The update re-sorts the table
# chronologically (2404 -> 2410), grouping the three workers together for
# each period, and moves the reduced "Valor Mora" (15600) from the first
# period block to the last one (2410).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Worker identity columns: Tipo Doc (B), N Doc (C), Nombre (D)
$workers = @(
    @("CC", "1044923835", "HENRY JOSE MORELO PINTO"),
    @("CC", "1143327560", "LUIS RAFAEL GUTIERREZ JULIO"),
    @("CC", "1143358068", "JOSE ANDRES SAENZ RUIZ")
)

$periods = @("2404", "2405", "2406", "2407", "2408", "2409", "2410")

$row = 16
foreach ($periodo in $periods) {
    foreach ($worker in $workers) {
        if ($periodo -eq "2410") {
            $valorMora = 15600
        } else {
            $valorMora = 52000
        }

        $ws.Cells.Item($row, 2).Value = $worker[0]
        $ws.Cells.Item($row, 3).Value = $worker[1]
        $ws.Cells.Item($row, 4).Value = $worker[2]
        $ws.Cells.Item($row, 5).Value = $periodo
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = 1300000

        $row = $row + 1
    }
}
